$d = $word.ActiveDocument

# Change "(Spring 2024)" -> "(Spring 2025)"
$d.Content.Find.Execute("(Spring 2024)", $true, $false, $false, $false, $false,
                         $true, 1, $false, "(Spring 2025)", 2) | Out-Null

# Change Study ID "2022 00000" -> "2024 00000"
$d.Content.Find.Execute("202200000", $true, $false, $false, $false, $false,
                         $true, 1, $false, "202400000", 2) | Out-Null
